$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 13.250951766967773
$ws.Range("P2").Value = 38.024566650390625
$ws.Range("Q2").Value = 29.036596298217773
$ws.Range("R2").Value = 8.9879703521728516
$ws.Range("S2").Value = 3.2160265445709229
$ws.Range("T2").Value = 35.548210144042969
$ws.Range("U2").Value = 8.6075477600097656
$ws.Range("V2").Value = 26.940662384033203
$ws.Range("W2").Value = 1
$ws.Range("X2").Value = 0
$ws.Range("Y2").Value = 17.446767807006836
$ws.Range("Z2").Value = 21.880966186523438
$ws.Range("AA2").Value = 1.1773288249969482
$ws.Range("AB2").Value = 20.703638076782227
$ws.Range("AC2").Value = 8.8961515426635742
$ws.Range("AD2").Value = 2.1978754997253418
$ws.Range("AE2").Value = 2.1978754997253418
$ws.Range("AF2").Value = 0
$ws.Range("AG2").Value = 1
$ws.Range("AI2").Value = 7.0325741767883301
$ws.Range("AJ2").Value = 29.231773376464844
$ws.Range("AK2").Value = 29.231773376464844
$ws.Range("AL2").Value = 0
$ws.Range("AM2").Value = 4.6432280540466309
$ws.Range("AN2").Value = 8.4582910537719727
$ws.Range("AO2").Value = 8.4582910537719727
$ws.Range("AP2").Value = 0
$ws.Range("AR2").Value = 0
$ws.Range("AS2").Value = 2.3537311553955078
$ws.Range("AT2").Value = 32.880050659179688
$ws.Range("AU2").Value = 9.2514591217041016
$ws.Range("AV2").Value = 23.628591537475586
$ws.Range("AW2").Value = 8.6124334335327148
$ws.Range("AX2").Value = 13.234278678894043
$ws.Range("AY2").Value = 13.234278678894043
$ws.Range("AZ2").Value = 0
$ws.Range("BC2").Value = 1.9467545747756958
$ws.Range("BD2").Value = 44.056102752685547
$ws.Range("BE2").Value = 39.308578491210938
$ws.Range("BF2").Value = 4.7475242614746094
$ws.Range("BG2").Value = 3.0712547302246094
$ws.Range("BH2").Value = 47.781539916992188
$ws.Range("BI2").Value = 43.099933624267578
$ws.Range("BJ2").Value = 4.6816062927246094
$ws.Range("BL2").Value = 0
$ws.Range("BM2").Value = 18.21574592590332
$ws.Range("BN2").Value = 42.9378662109375
$ws.Range("BO2").Value = 7.8077750205993652
$ws.Range("BP2").Value = 35.130092620849609
$ws.Range("BQ2").Value = 14.384234428405762
$ws.Range("BR2").Value = 46.834102630615234
$ws.Range("BS2").Value = 3.9423618316650391
$ws.Range("BT2").Value = 42.891738891601563
$ws.Range("BU2").Value = 10.666536331176758
$ws.Range("BV2").Value = 7.4639077186584473
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 13.250951766967773
$ws.Range("P3").Value = 38.024566650390625
$ws.Range("Q3").Value = 29.036596298217773
$ws.Range("R3").Value = 8.9879703521728516
$ws.Range("S3").Value = 3.2160265445709229
$ws.Range("T3").Value = 35.548210144042969
$ws.Range("U3").Value = 8.6075477600097656
$ws.Range("V3").Value = 26.940662384033203
$ws.Range("W3").Value = 1
$ws.Range("X3").Value = 0
$ws.Range("Y3").Value = 17.446767807006836
$ws.Range("Z3").Value = 21.880966186523438
$ws.Range("AA3").Value = 1.1773288249969482
$ws.Range("AB3").Value = 20.703638076782227
$ws.Range("AC3").Value = 8.8961515426635742
$ws.Range("AD3").Value = 2.1978754997253418
$ws.Range("AE3").Value = 2.1978754997253418
$ws.Range("AF3").Value = 0
$ws.Range("AG3").Value = 1
$ws.Range("AI3").Value = 7.0325741767883301
$ws.Range("AJ3").Value = 29.231773376464844
$ws.Range("AK3").Value = 29.231773376464844
$ws.Range("AL3").Value = 0
$ws.Range("AM3").Value = 4.6432280540466309
$ws.Range("AN3").Value = 8.4582910537719727
$ws.Range("AO3").Value = 8.4582910537719727
$ws.Range("AP3").Value = 0
$ws.Range("AR3").Value = 0
$ws.Range("AS3").Value = 2.3537311553955078
$ws.Range("AT3").Value = 32.880050659179688
$ws.Range("AU3").Value = 9.2514591217041016
$ws.Range("AV3").Value = 23.628591537475586
$ws.Range("AW3").Value = 8.6124334335327148
$ws.Range("AX3").Value = 13.234278678894043
$ws.Range("AY3").Value = 13.234278678894043
$ws.Range("AZ3").Value = 0
$ws.Range("BC3").Value = 1.9467545747756958
$ws.Range("BD3").Value = 44.056102752685547
$ws.Range("BE3").Value = 39.308578491210938
$ws.Range("BF3").Value = 4.7475242614746094
$ws.Range("BG3").Value = 3.0712547302246094
$ws.Range("BH3").Value = 47.781539916992188
$ws.Range("BI3").Value = 43.099933624267578
$ws.Range("BJ3").Value = 4.6816062927246094
$ws.Range("BL3").Value = 0
$ws.Range("BM3").Value = 18.21574592590332
$ws.Range("BN3").Value = 42.9378662109375
$ws.Range("BO3").Value = 7.8077750205993652
$ws.Range("BP3").Value = 35.130092620849609
$ws.Range("BQ3").Value = 14.384234428405762
$ws.Range("BR3").Value = 46.834102630615234
$ws.Range("BS3").Value = 3.9423618316650391
$ws.Range("BT3").Value = 42.891738891601563
$ws.Range("BU3").Value = 10.666536331176758
$ws.Range("BV3").Value = 7.4639077186584473
